# Apply the changes described by the diff:
# 1. Update the "Dt. Referencia" column (G2:G274) from 45615 (2024-11-19) to 45617 (2024-11-21)
# 2. Correct the balance values for two rows (51 and 118)
# 3. Rename the worksheet to reflect the new export timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update date column G for all data rows (2 through 274)
$ws.Range("G2:G274").Value = 45617

# 2. Fix the mismatched balance values
$ws.Range("E51").Value = 25411.56
$ws.Range("H51").Value = 25411.56

$ws.Range("E118").Value = 9717.06
$ws.Range("H118").Value = 9717.06

# 3. Rename the worksheet (and update the workbook's sheet name) to the new export id
$ws.Name = "IClientBalance-20241121-094552-"
